# fixed negative comment scores
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 3: post bd0ww1
$ws.Range("B3").Value = 0.00186046511627907
$ws.Range("F3").Value = 49

# Row 4: post blexov
$ws.Range("B4").Value = 0.9448717948717948
$ws.Range("G4").Value = 8

# Row 5: post aos6vn
$ws.Range("B5").Value = 0.01162790697674419
$ws.Range("F5").Value = 7

# Row 6: post b1cbcz
$ws.Range("F6").Value = 7

# Row 10: post cjudzm
$ws.Range("B10").Value = 0.9850746268656716
$ws.Range("D10").Value = 655
$ws.Range("G10").Value = 5

# Row 17: post hql2q4
$ws.Range("D17").Value = 582
$ws.Range("G17").Value = 19
